$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 176, shifting existing rows 176:215 down to 177:216
$ws.Rows("176:176").Insert()

# Populate the newly inserted row 176 with the new weekly record
$ws.Cells.Item(176, 1).Value = 4
$ws.Cells.Item(176, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(176, 3).Value = "Los Lagos"
$ws.Cells.Item(176, 4).Value = 44641
$ws.Cells.Item(176, 5).Value = 10
$ws.Cells.Item(176, 6).Value = "Fruta"
$ws.Cells.Item(176, 7).Value = 100108
$ws.Cells.Item(176, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(176, 9).Value = 100108005
$ws.Cells.Item(176, 10).Value = "Piña"
$ws.Cells.Item(176, 11).Value = "Caramelo"
$ws.Cells.Item(176, 12).Value = "Tercera"
$ws.Cells.Item(176, 13).Value = 50
$ws.Cells.Item(176, 14).Value = 17000
$ws.Cells.Item(176, 15).Value = 18000
$ws.Cells.Item(176, 16).Value = 17500
$ws.Cells.Item(176, 17).Value = "$/caja 16 unidades"
$ws.Cells.Item(176, 18).Value = "Ecuador"
$ws.Cells.Item(176, 19).Value = 1094
$ws.Cells.Item(176, 20).Value = 16
